$d = $word.ActiveDocument

# 1) "亲爱的朋友，" -> "Dear friends."
$found0 = $d.Content.Find.Execute(
    "亲爱的朋友，", $false, $false, $false, $false, $false, $true, 1, $false,
    "Dear friends.", 2)
if (-not $found0) {
    throw "Could not find greeting '亲爱的朋友，' to replace"
}

# 2) Translate the "抱歉，我们用了将近一个月的时间，才把 2015 年第一个季度的工作
#    报告呈现在您的面前。" sentence to English while keeping the remaining
#    Chinese text of that paragraph untouched. The lookups are chained through
#    progressively narrower Range objects so that only the first occurrence
#    (the one immediately following the greeting) is touched, since "2015"
#    also appears later in the document.
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "抱歉，我们用了将近一个月的时间，才把", $false, $false, $false, $false, $false, $true, 1, $false,
    "Sorry, it took us nearly a month to present the work report for the first quarter of ", 2)
if (-not $found1) {
    throw "Could not find '抱歉，我们用了将近一个月的时间，才把' to replace"
}
$r1.Collapse(0)

$r2 = $d.Range($r1.End, $r1.End + 300)
$found2 = $r2.Find.Execute(
    "2015", $false, $false, $false, $false, $false, $true, 1, $false,
    "2015 ", 2)
if (-not $found2) {
    throw "Could not find '2015' following the greeting to replace"
}
$r2.Collapse(0)

$r3 = $d.Range($r2.End, $r2.End + 300)
$found3 = $r3.Find.Execute(
    "年第一个季度的工作报告呈现在您的面前。 ", $false, $false, $false, $false, $false, $true, 1, $false,
    "to you. ", 2)
if (-not $found3) {
    throw "Could not find '年第一个季度的工作报告呈现在您的面前。 ' to replace"
}
